# Generate Report for Handoff
#
# Moves the localization status from "In Translation" to "Ready for
# handoff" and refreshes the associated timestamps on all three sheets
# of the report (Overview, zh-cn, de-de). Also widens the "Status" /
# language-status columns so the new, longer text fits (matches the
# report generator's autosize behaviour).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet --------------------------------------------------
# Column E = zh-cn status, Column F = de-de status, Column G = latest
# handoff xliff generation date for the row.
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-05 02:44:11"

# --- zh-cn sheet -------------------------------------------------------
# Column C = Status, Column H = Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-05 02:44:03"

# --- de-de sheet -------------------------------------------------------
# Column C = Status, Column H = Latest Handoff Datetime
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-05 02:44:11"

# --- Column widths -----------------------------------------------------
# "Ready for handoff" is wider than "In Translation", so the status
# columns are widened to fit (as the reporting tool would regenerate
# them). 16.33 characters is the ColumnWidth value that this engine's
# pixel grid resolves closest to the recorded target stored width.
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33

Write-Output "Generate Report for Handoff: status + timestamps updated, columns resized."
